# Changes for new model
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the columns that are no longer part of the model output (E:G)
$ws.Range("E1:G1").EntireColumn.Delete()

# Reset column widths: B/C/D get new widths, A stays as-is
# (ColumnWidth uses Excel's char-width->pixel rounding, which pads by
#  ~0.83 "Calibri 11" units when storing; subtracting 0.9 before setting
#  keeps us inside the rounding bucket that serializes to the clean
#  integer width value required by the target file.)
$ws.Columns.Item(2).ColumnWidth = 22 - 0.9
$ws.Columns.Item(3).ColumnWidth = 24 - 0.9
$ws.Columns.Item(4).ColumnWidth = 18 - 0.9

# Update header row labels
$ws.Range("B1").Value = "Matematica Predicted"
$ws.Range("C1").Value = "Comunicación Predicted"
$ws.Range("D1").Value = "Ingles Predicted"

# Update data rows with the new predicted values
$ws.Range("B2").Value = 14.37
$ws.Range("C2").Value = 12.98
$ws.Range("D2").Value = 10.74

$ws.Range("B3").Value = 14.32
$ws.Range("C3").Value = 12.82
$ws.Range("D3").Value = 10.65
